# This script updates the "想去人数" (want-to-go count) values in the
# "展览" and "全部类型" worksheets, incrementing a set of specific cells
# to reflect newly generated data (per commit message: "Update gh-pages
# to output generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# Map of worksheet name -> cell address -> new value
$updates = @{
    "展览" = @{
        "F3"  = 281
        "F4"  = 1215
        "F10" = 3418
        "F11" = 122
        "F12" = 82
        "F16" = 583
        "F18" = 706
        "F22" = 53
        "F23" = 62
        "F24" = 2517
        "F25" = 5033
        "F31" = 2215
        "F36" = 162
        "F38" = 457
        "F39" = 782
        "F43" = 461
    }
    "全部类型" = @{
        "F3"  = 281
        "F4"  = 1215
        "F10" = 3418
        "F11" = 122
        "F12" = 82
        "F17" = 583
        "F19" = 706
        "F23" = 53
        "F24" = 62
        "F25" = 2517
        "F26" = 5033
        "F32" = 2215
        "F37" = 162
        "F39" = 457
        "F40" = 782
        "F44" = 461
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellUpdates = $updates[$sheetName]
    foreach ($addr in $cellUpdates.Keys) {
        $ws.Range($addr).Value = $cellUpdates[$addr]
    }
}
